# Minor refactoring to SpreadsheetImporter
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Machray" -> "Machray Hall" (the shared string used by rows 6-10, column A)
$ws.Range("A6:A10").Value = "Machray Hall"

# Give the sheet tab a white color
$ws.Tab.Color = 16777215

# Move the active selection to A11
$ws.Range("A11").Select()
